$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7689787935149184
$ws.Range("D2").Value = 0.09934938549012529
$ws.Range("E2").Value = 0.1193460124957659
$ws.Range("F2").Value = 2.006493858668705
$ws.Range("G2").Value = 1.370177519816579
$ws.Range("H2").Value = 1.279577249586623
$ws.Range("I2").Value = 1.275467134579213
$ws.Range("J2").Value = 0.1660899761622581
$ws.Range("K2").Value = 0.6920889448149978
$ws.Range("L2").Value = 0.3417046917287081
$ws.Range("N2").Value = 2.250971858947203
$ws.Range("B3").Value = 0.7457963122945159
$ws.Range("D3").Value = 0.09915605637133851
$ws.Range("E3").Value = 0.1196920202080007
$ws.Range("F3").Value = 2.008646777177354
$ws.Range("G3").Value = 1.368971464488354
$ws.Range("H3").Value = 1.284041327579843
$ws.Range("I3").Value = 1.284220363916905
$ws.Range("J3").Value = 0.1665825441896436
$ws.Range("K3").Value = 0.6340528164168973
$ws.Range("L3").Value = 0.3307997776868206
$ws.Range("N3").Value = 2.272235882313568
$ws.Range("B4").Value = 0.7318758648983703
$ws.Range("D4").Value = 0.09905504272061449
$ws.Range("E4").Value = 0.1199213408156465
$ws.Range("F4").Value = 2.010911322070747
$ws.Range("G4").Value = 1.368986905682178
$ws.Range("H4").Value = 1.287313424717141
$ws.Range("I4").Value = 1.290083090010299
$ws.Range("J4").Value = 0.1669040084987445
$ws.Range("K4").Value = 0.5986448454915774
$ws.Range("L4").Value = 0.3242642547458559
$ws.Range("N4").Value = 2.285954355169864
$ws.Range("B5").Value = 0.7262825694539572
$ws.Range("D5").Value = 0.09901835041024043
$ws.Range("E5").Value = 0.1200190439563547
$ws.Range("F5").Value = 2.012071311825395
$ws.Range("G5").Value = 1.369183239710594
$ws.Range("H5").Value = 1.28878048304658
$ws.Range("I5").Value = 1.292595034939247
$ws.Range("J5").Value = 0.1670398027595885
$ws.Range("K5").Value = 0.5842732512741406
$ws.Range("L5").Value = 0.3216413758860455
$ws.Range("N5").Value = 2.291711294995345
$ws.Range("B6").Value = 0.7253586178865135
$ws.Range("D6").Value = 0.09901252843044972
$ws.Range("E6").Value = 0.1200355246890576
$ws.Range("F6").Value = 2.0122782556613
$ws.Range("G6").Value = 1.369227317402562
$ws.Range("H6").Value = 1.28903216201401
$ws.Range("I6").Value = 1.293019562647846
$ws.Range("J6").Value = 0.1670626412148501
$ws.Range("K6").Value = 0.5818903430918283
$ws.Range("L6").Value = 0.3212082927910558
$ws.Range("N6").Value = 2.292677290239206
$ws.Range("B7").Value = 0.731800109582764
$ws.Range("D7").Value = 0.09905452973854167
$ws.Range("E7").Value = 0.119922641238658
$ws.Range("F7").Value = 2.01092600563377
$ws.Range("G7").Value = 1.368988784110599
$ws.Range("H7").Value = 1.287332668704209
$ws.Range("I7").Value = 1.290116469523007
$ws.Range("J7").Value = 0.1669058204386258
$ws.Range("K7").Value = 0.5984507919046678
$ws.Range("L7").Value = 0.32422871793176
$ws.Range("N7").Value = 2.286031320897616
$ws.Range("B8").Value = 0.7609206706120801
$ws.Range("D8").Value = 0.09927906736673364
$ws.Range("E8").Value = 0.1194618210230574
$ws.Range("F8").Value = 2.007040614746728
$ws.Range("G8").Value = 1.369604736042092
$ws.Range("H8").Value = 1.281006281121549
$ws.Range("I8").Value = 1.278383963894637
$ws.Range("J8").Value = 0.166255871376638
$ws.Range("K8").Value = 0.6720313694960964
$ws.Range("L8").Value = 0.337911515418412
$ws.Range("N8").Value = 2.258166250227326
$ws.Range("B9").Value = 0.8204957013398371
$ws.Range("D9").Value = 0.0998588184785234
$ws.Range("E9").Value = 0.1186915337197059
$ws.Range("F9").Value = 2.006897042123853
$ws.Range("G9").Value = 1.376815205956561
$ws.Range("H9").Value = 1.272811233569783
$ws.Range("I9").Value = 1.259247209835006
$ws.Range("J9").Value = 0.1651317934355796
$ws.Range("K9").Value = 0.8181051498983152
$ws.Range("L9").Value = 0.3660100126754742
$ws.Range("N9").Value = 2.208776246321495
$ws.Range("B10").Value = 0.8657503921433545
$ws.Range("D10").Value = 0.1003686072180585
$ws.Range("E10").Value = 0.1182062696128532
$ws.Range("F10").Value = 2.011346922467013
$ws.Range("G10").Value = 1.385781028091628
$ws.Range("H10").Value = 1.269354134285976
$ws.Range("I10").Value = 1.247543181412105
$ws.Range("J10").Value = 0.1643969861105212
$ws.Range("K10").Value = 0.9265044069567807
$ws.Range("L10").Value = 0.3874232299922511
$ws.Range("N10").Value = 2.175688629203043
$ws.Range("B11").Value = 0.8866563052121421
$ws.Range("D11").Value = 0.1006184940121102
$ws.Range("E11").Value = 0.1180028910264426
$ws.Range("F11").Value = 2.014360259595577
$ws.Range("G11").Value = 1.390658481069892
$ws.Range("H11").Value = 1.268337519063763
$ws.Range("I11").Value = 1.242729456239147
$ws.Range("J11").Value = 0.1640823311729349
$ws.Range("K11").Value = 0.9760513440485283
$ws.Range("L11").Value = 0.39733124515989
$ws.Range("N11").Value = 2.161330189369201
$ws.Range("B12").Value = 0.8946182961465752
$ws.Range("D12").Value = 0.1007156841976595
$ws.Range("E12").Value = 0.117928364047212
$ws.Range("F12").Value = 2.015643490701521
$ws.Range("G12").Value = 1.392620440937378
$ws.Range("H12").Value = 1.268032446665501
$ws.Range("I12").Value = 1.240979973546608
$ws.Range("J12").Value = 0.1639659889568525
$ws.Range("K12").Value = 0.9948470123673587
$ws.Range("L12").Value = 0.4011070780861985
$ws.Range("N12").Value = 2.155992708346812
$ws.Range("B13").Value = 0.8929015295031206
$ws.Range("D13").Value = 0.1006946388927261
$ws.Range("E13").Value = 0.1179443042624495
$ws.Range("F13").Value = 2.015360803223686
$ws.Range("G13").Value = 1.392192782668587
$ws.Range("H13").Value = 1.268094596929458
$ws.Range("I13").Value = 1.241353493370099
$ws.Range("J13").Value = 0.1639909204745811
$ws.Range("K13").Value = 0.9907975559030149
$ws.Range("L13").Value = 0.4002928248982727
$ws.Range("N13").Value = 2.157137793126019
$ws.Range("B14").Value = 0.8873104367184226
$ws.Range("D14").Value = 0.1006264386575104
$ws.Range("E14").Value = 0.1179967098378025
$ws.Range("F14").Value = 2.014462983436417
$ws.Range("G14").Value = 1.390817587859402
$ws.Range("H14").Value = 1.268310819642863
$ws.Range("I14").Value = 1.242584055005594
$ws.Range("J14").Value = 0.1640727033507252
$ws.Range("K14").Value = 0.9775970104263081
$ws.Range("L14").Value = 0.3976414073472938
$ws.Range("N14").Value = 2.160889072156039
$ws.Range("B15").Value = 0.8838916240444519
$ws.Range("D15").Value = 0.1005849972297987
$ws.Range("E15").Value = 0.1180291334986481
$ws.Range("F15").Value = 2.013931552090668
$ws.Range("G15").Value = 1.389990217347361
$ws.Range("H15").Value = 1.268453665587742
$ws.Range("I15").Value = 1.243347363432726
$ws.Range("J15").Value = 0.1641231634900353
$ws.Range("K15").Value = 0.9695156177118918
$ws.Range("L15").Value = 0.3960204437019144
$ws.Range("N15").Value = 2.163199830776895
$ws.Range("B16").Value = 0.8643905207594003
$ws.Range("D16").Value = 0.1003526364550282
$ws.Range("E16").Value = 0.1182199095823009
$ws.Range("F16").Value = 2.011169889395589
$ws.Range("G16").Value = 1.385478358606321
$ws.Range("H16").Value = 1.269431756233516
$ws.Range("I16").Value = 1.247868038878956
$ws.Range("J16").Value = 0.1644179433933224
$ws.Range("K16").Value = 0.9232710908352431
$ws.Range("L16").Value = 0.3867790673728422
$ws.Range("N16").Value = 2.17664094228947
$ws.Range("B17").Value = 0.8525086259201657
$ws.Range("D17").Value = 0.100214680221832
$ws.Range("E17").Value = 0.1183413865393468
$ws.Range("F17").Value = 2.009728964271289
$ws.Range("G17").Value = 1.382915156347551
$ws.Range("H17").Value = 1.270174154673967
$ws.Range("I17").Value = 1.250772049728397
$ws.Range("J17").Value = 0.1646037977916821
$ws.Range("K17").Value = 0.8949615110918501
$ws.Range("L17").Value = 0.3811524666188006
$ws.Range("N17").Value = 2.185064233538087
$ws.Range("B18").Value = 0.8457045632981135
$ws.Range("D18").Value = 0.100137025250433
$ws.Range("E18").Value = 0.1184128924440531
$ws.Range("F18").Value = 2.008993282394954
$ws.Range("G18").Value = 1.381516056256103
$ws.Range("H18").Value = 1.270653502900046
$ws.Range("I18").Value = 1.252490412489031
$ws.Range("J18").Value = 0.1647125429189433
$ws.Range("K18").Value = 0.8787007959302287
$ws.Range("L18").Value = 0.3779319249970143
$ws.Range("N18").Value = 2.189974339434833
$ws.Range("B19").Value = 0.843406008998727
$ws.Range("D19").Value = 0.1001110242041108
$ws.Range("E19").Value = 0.1184373843250022
$ws.Range("F19").Value = 2.008760185666603
$ws.Range("G19").Value = 1.381055255252861
$ws.Range("H19").Value = 1.270824793307526
$ws.Range("I19").Value = 1.253080474840935
$ws.Range("J19").Value = 0.1647496795948498
$ws.Range("K19").Value = 0.8731990280874697
$ws.Range("L19").Value = 0.3768442108099208
$ws.Range("N19").Value = 2.191648022319781
$ws.Range("B20").Value = 0.8537703637485947
$ws.Range("D20").Value = 0.1002291907585331
$ws.Range("E20").Value = 0.1183282859057631
$ws.Range("F20").Value = 2.009872718277862
$ws.Range("G20").Value = 1.383180231405021
$ws.Range("H20").Value = 1.270089708502454
$ws.Range("I20").Value = 1.25045793983643
$ws.Range("J20").Value = 0.1645838222276028
$ws.Range("K20").Value = 0.8979728199412307
$ws.Range("L20").Value = 0.3817498008383353
$ws.Range("N20").Value = 2.184160806119461
$ws.Range("B21").Value = 0.888951449051973
$ws.Range("D21").Value = 0.100646401323722
$ws.Range("E21").Value = 0.1179812496016841
$ws.Range("F21").Value = 2.014722837713464
$ws.Range("G21").Value = 1.391218394743078
$ws.Range("H21").Value = 1.268245141859111
$ws.Range("I21").Value = 1.242220618268057
$ws.Range("J21").Value = 0.1640486055288073
$ws.Range("K21").Value = 0.981473430796683
$ws.Range("L21").Value = 0.3984195458753561
$ws.Range("N21").Value = 2.159784522969176
$ws.Range("B22").Value = 0.9122085528301511
$ws.Range("D22").Value = 0.1009340043183826
$ws.Range("E22").Value = 0.1177689392715129
$ws.Range("F22").Value = 2.01872116862053
$ws.Range("G22").Value = 1.397142003933467
$ws.Range("H22").Value = 1.267505289109437
$ws.Range("I22").Value = 1.237264695075382
$ws.Range("J22").Value = 0.1637151895185003
$ws.Range("K22").Value = 1.036239854863993
$ws.Range("L22").Value = 0.4094533306887058
$ws.Range("N22").Value = 2.144434651543421
$ws.Range("B23").Value = 0.8997718010666063
$ws.Range("D23").Value = 0.1007791462113161
$ws.Range("E23").Value = 0.1178809299297958
$ws.Range("F23").Value = 2.016511402932821
$ws.Range("G23").Value = 1.393919105332415
$ws.Range("H23").Value = 1.267857572288449
$ws.Range("I23").Value = 1.239870647623775
$ws.Range("J23").Value = 0.1638916443442211
$ws.Range("K23").Value = 1.006992421198191
$ws.Range("L23").Value = 0.403551705668562
$ws.Range("N23").Value = 2.152573946001789
$ws.Range("B24").Value = 0.8531998479423919
$ws.Range("D24").Value = 0.1002226253798391
$ws.Range("E24").Value = 0.1183342035090598
$ws.Range("F24").Value = 2.009807438247975
$ws.Range("G24").Value = 1.383060158853837
$ws.Range("H24").Value = 1.270127722968368
$ws.Range("I24").Value = 1.250599796878376
$ws.Range("J24").Value = 0.164592847277429
$ws.Range("K24").Value = 0.8966113609472188
$ws.Range("L24").Value = 0.3814797015907203
$ws.Range("N24").Value = 2.184569035569456
$ws.Range("B25").Value = 0.8041168202159099
$ws.Range("D25").Value = 0.09968717352927214
$ws.Range("E25").Value = 0.1188857047724921
$ws.Range("F25").Value = 2.006135772623267
$ws.Range("G25").Value = 1.374221041898522
$ws.Range("H25").Value = 1.274577744961448
$ws.Range("I25").Value = 1.264010236293686
$ws.Range("J25").Value = 0.1654198478637481
$ws.Range("K25").Value = 0.7783983018283323
$ws.Range("L25").Value = 0.3582733285847581
$ws.Range("N25").Value = 2.221575506376224

Write-Host "Applied 264 cell updates"
